# Applies the "added mango and shaga testnet" edit to AIRDROP.xlsx
# - Clears the stray "NO" header value from A1 (header row loses col-A label)
# - Fills in the three previously-blank rows (43-45) with new project rows:
#     43: SHAGA GLOB
#     44: MANGO NETWORK  (link becomes a clickable hyperlink, like the rest of column C)
#     45: SINGULARITY FINANCE
# - Restores the selection to C52 to mirror the saved cursor position

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row: A1 previously held the shared string "NO" - clear it out.
$ws.Range("A1").Value = ""

# Row 43 - SHAGA GLOB
$ws.Range("B43").Value = "SHAGA GLOB"
$ws.Range("C43").Value = "https://glob.shaga.xyz/?r=Bdgdwjwdtb"
$ws.Range("D43").Value = "CHAIN/DEPIN"
$ws.Range("E43").Value = "PROCESSING"

# Row 44 - MANGO NETWORK (link pasted in before the name, like the source edit)
$ws.Range("C44").Value = "https://task.testnet.mangonetwork.io/?invite=FHKOiL"
$ws.Range("B44").Value = "MANGO NETWORK"
$ws.Range("D44").Value = "TESTNET"
$ws.Range("E44").Value = "PROCESSING"
$ws.Hyperlinks.Add($ws.Range("C44"), "https://task.testnet.mangonetwork.io/?invite=FHKOiL")
# Hyperlinks.Add stamps the built-in blue/underline "Hyperlink" style; every
# other linked cell in column C instead keeps the sheet's plain link style
# (style index 10, seen on C2/C3/...), so copy that formatting back over.
$ws.Range("C2").Copy()
$ws.Range("C44").PasteSpecial(-4122)  # xlPasteFormats

# Row 45 - SINGULARITY FINANCE (link pasted in before the name, like the source edit)
$ws.Range("C45").Value = "https://t.me/KeoAirDropFreeNe/323/37670"
$ws.Range("B45").Value = "SINGULARITY FINANCE"
$ws.Range("D45").Value = "FINANCE/TESTNET"
$ws.Range("E45").Value = "PROCESSING"

# Match the saved cursor/selection state recorded in the workbook.
$ws.Range("C52").Select()
